# Refresh "top-crypto" data (query refresh as of Fri 12-31-2021) and
# adjust column D width (which also changes the current selection to the
# full column, matching Excel's behaviour when a column is resized via
# its header boundary).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (BTC) ---
$ws.Range("E2").Value = 47303.409690733315
$ws.Range("F2").Value = 40.342199999999998
$ws.Range("H2").Value = 18915200
$ws.Range("I2").Value = 28067070886.93898
$ws.Range("J2").Value = 894753454982.15881

# --- Row 3 (ETH) ---
$ws.Range("E3").Value = 3728.2190797676931
$ws.Range("F3").Value = 19.996200000000002
$ws.Range("H3").Value = 118956833.874
$ws.Range("I3").Value = 14453683338.12944
$ws.Range("J3").Value = 443497137717.80261

# --- Row 4 (BNB) ---
$ws.Range("E4").Value = 517.18007448864648
$ws.Range("F4").Value = 3.8965999999999998
$ws.Range("I4").Value = 2065289968.3729904
$ws.Range("J4").Value = 86266230147.431747

# --- Row 5 (USDT) ---
$ws.Range("E5").Value = 1.0018538229298697
$ws.Range("F5").Value = 3.5348000000000002
$ws.Range("I5").Value = 57613421140.380203
$ws.Range("J5").Value = 78490120046.346497

# --- Row 6 (SOL) --- row height shrinks (autofit for new column width)
$ws.Rows("6").RowHeight = 57.6
$ws.Range("E6").Value = 173.24823906902782
$ws.Range("F6").Value = 2.4218999999999999
$ws.Range("H6").Value = 309484316.31080955
$ws.Range("I6").Value = 1722954113.0950823
$ws.Range("J6").Value = 53617612820.329758

# --- Row 7 (ADA) ---
$ws.Range("E7").Value = 1.345814886076151
$ws.Range("F7").Value = 2.0356000000000001
$ws.Range("H7").Value = 33485576474.460999
$ws.Range("I7").Value = 1343839034.0396678
$ws.Range("J7").Value = 45065387288.170967

# --- Row 8 (USDC) ---
$ws.Range("E8").Value = 1.0012429923053927
$ws.Range("F8").Value = 1.9111
$ws.Range("H8").Value = 42281416428.680771
$ws.Range("I8").Value = 3323193156.9762902
$ws.Range("J8").Value = 42333971903.962723

# --- Row 9 (XRP) ---
$ws.Range("E9").Value = 0.83979062876816035
$ws.Range("F9").Value = 1.8021
$ws.Range("I9").Value = 1994764321.754813
$ws.Range("J9").Value = 39920257493.881599

# --- Row 10 (LUNA) ---
$ws.Range("E10").Value = 83.895118345524111
$ws.Range("F10").Value = 1.3694999999999999
$ws.Range("H10").Value = 362053961.31293142
$ws.Range("I10").Value = 2276388333.7667866
$ws.Range("J10").Value = 30374559931.81419

# --- Row 11 (DOT) --- row height shrinks (autofit for new column width)
$ws.Rows("11").RowHeight = 57.6
$ws.Range("E11").Value = 27.546714612609406
$ws.Range("F11").Value = 1.2287999999999999
$ws.Range("I11").Value = 1403985685.0398903
$ws.Range("J11").Value = 27204565546.439121

# --- Row 12 (AVAX) ---
$ws.Range("E12").Value = 101.73949722896118
$ws.Range("F12").Value = 1.1172
$ws.Range("H12").Value = 243245553.82933065
$ws.Range("I12").Value = 1043478184.0910224
$ws.Range("J12").Value = 24747680349.776314

# --- Row 13 (DOGE) ---
$ws.Range("E13").Value = 0.17145150169927939
$ws.Range("F13").Value = 1.0274000000000001
$ws.Range("I13").Value = 685970736.1827482
$ws.Range("J13").Value = 22746601770.807987

# --- Row 14 (SHIB) ---
$ws.Range("E14").Value = 0.000033755709583294796
$ws.Range("F14").Value = 0.83720000000000006
$ws.Range("I14").Value = 1099112032.1563272
$ws.Range("J14").Value = 18534020584.600048

# --- Row 15 (MATIC) --- row height shrinks (autofit for new column width)
$ws.Rows("15").RowHeight = 57.6
$ws.Range("E15").Value = 2.5613030304846336
$ws.Range("F15").Value = 0.82869999999999999
$ws.Range("I15").Value = 2049691143.1838143
$ws.Range("J15").Value = 18346338019.530811

# --- Row 16 (BUSD) ---
$ws.Range("E16").Value = 1.0015920886635199
$ws.Range("F16").Value = 0.66200000000000003
$ws.Range("I16").Value = 3917340742.6773281
$ws.Range("J16").Value = 14656494194.31196

# --- Row 17 (CRO) ---
$ws.Range("E17").Value = 0.55845310652645797
$ws.Range("F17").Value = 0.63729999999999998
$ws.Range("I17").Value = 216459770.07655615
$ws.Range("J17").Value = 14108208476.517839

# --- Row 18 (WBTC) ---
$ws.Range("E18").Value = 47721.472895008992
$ws.Range("F18").Value = 0.55640000000000001
$ws.Range("I18").Value = 244192154.78925464
$ws.Range("J18").Value = 12318865860.694613

# --- Row 19: was UNI/Uniswap, now ALGO/Algorand (ranking swap with row 20) ---
$ws.Range("A19").Value = "ALGO"
$ws.Range("B19").Value = "Algorand"
$ws.Range("E19").Value = 1.7159673112700276
$ws.Range("F19").Value = 0.49349999999999999
$ws.Range("G19").Value = 10000000000
$ws.Range("H19").Value = 6366823990.8037624
$ws.Range("I19").Value = 900407624.33687866
$ws.Range("J19").Value = 10925261844.829041

# --- Row 20: was ALGO/Algorand, now UNI/Uniswap (ranking swap with row 19) ---
$ws.Range("A20").Value = "UNI"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("E20").Value = 17.366878960125263
$ws.Range("F20").Value = 0.49209999999999998
$ws.Range("G20").Value = 1000000000
$ws.Range("H20").Value = 627264509.89428568
$ws.Range("I20").Value = 311734607.4073956
$ws.Range("J20").Value = 10893626819.316357

# --- Row 21 (LTC) ---
$ws.Range("E21").Value = 148.02146847132292
$ws.Range("F21").Value = 0.46250000000000002
$ws.Range("H21").Value = 69305944.562391832
$ws.Range("I21").Value = 917926295.66619599
$ws.Range("J21").Value = 10258767687.917336

# --- Row 22 (LINK) --- row height shrinks (autofit for new column width)
$ws.Rows("22").RowHeight = 57.6
$ws.Range("E22").Value = 1.0018887076383778
$ws.Range("F22").Value = 0.45600000000000002
$ws.Range("H22").Value = 10094236068.464312
$ws.Range("I22").Value = 121476175.33419564
$ws.Range("J22").Value = 10113301129.23041

# --- Row 23 (UST) ---
$ws.Range("C23").Value = 23
$ws.Range("E23").Value = 19.801557746047585
$ws.Range("F23").Value = 0.41770000000000002
$ws.Range("I23").Value = 990302101.56686223
$ws.Range("J23").Value = 9247516580.9887047

# --- Column D width change (49.109375 -> 57.21875) ---
$ws.Columns("D").ColumnWidth = 57.21875

# --- Selection becomes the whole column D (as happens when resizing a
#     column by dragging its header border) ---
$ws.Range("A1:XFD1048576").Select
